$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear data cells (B:L) for rows that become section-header-only rows
$ws.Range("B8:L8").ClearContents()
$ws.Range("B14:L14").ClearContents()
$ws.Range("B20:L20").ClearContents()
$ws.Range("B26:L26").ClearContents()
$ws.Range("B30:L30").ClearContents()

# Row 7
$ws.Range("B7,C7,D7,E7,F7,G7,H7,I7,J7,K7,L7").NumberFormat = "@"
$ws.Range("A7").Value = "Less than Primary"
$ws.Range("B7").Value = "0.44"
$ws.Range("C7").Value = "0.19"
$ws.Range("D7").Value = "0.32"
$ws.Range("E7").Value = "0.12"
$ws.Range("F7").Value = "0.25"
$ws.Range("G7").Value = "0.11"
$ws.Range("H7").Value = "0.11"
$ws.Range("I7").Value = "0.12"
$ws.Range("J7").Value = "0.06"
$ws.Range("K7").Value = "0.01"
$ws.Range("L7").Value = "0.04"

# Row 8
$ws.Range("A8").Value = "Education Completed"

# Row 9
$ws.Range("B9,C9,D9,E9,F9,G9,H9,I9,J9,L9").NumberFormat = "@"
$ws.Range("B9").Value = "0.31"
$ws.Range("C9").Value = "0.3"
$ws.Range("D9").Value = "0.31"
$ws.Range("E9").Value = "0.26"
$ws.Range("F9").Value = "0.25"
$ws.Range("G9").Value = "0.16"
$ws.Range("H9").Value = "0.13"
$ws.Range("I9").Value = "0.24"
$ws.Range("J9").Value = "0.25"
$ws.Range("L9").Value = "0.15"

# Row 10
$ws.Range("A10").Value = "Secondary"

# Row 12
$ws.Range("B12,C12,D12,E12,F12,G12,H12,I12,J12,K12,L12").NumberFormat = "@"
$ws.Range("A12").Value = "Household Size"
$ws.Range("B12").Value = "3.55"
$ws.Range("C12").Value = "2.35"
$ws.Range("D12").Value = "3.07"
$ws.Range("E12").Value = "2.48"
$ws.Range("F12").Value = "3.21"
$ws.Range("G12").Value = "2.83"
$ws.Range("H12").Value = "2.58"
$ws.Range("I12").Value = "2.39"
$ws.Range("J12").Value = "2.13"
$ws.Range("K12").Value = "1.91"
$ws.Range("L12").Value = "2.25"

# Row 13
$ws.Range("B13,C13,D13,E13,F13,G13,H13,I13,J13,K13,L13").NumberFormat = "@"
$ws.Range("A13").Value = "Lives Alone"
$ws.Range("B13").Value = "0.15"
$ws.Range("C13").Value = "0.3"
$ws.Range("D13").Value = "0.23"
$ws.Range("E13").Value = "0.25"
$ws.Range("F13").Value = "0.19"
$ws.Range("G13").Value = "0.19"
$ws.Range("H13").Value = "0.23"
$ws.Range("I13").Value = "0.26"
$ws.Range("J13").Value = "0.37"
$ws.Range("K13").Value = "0.3"
$ws.Range("L13").Value = "0.3"

# Row 14
$ws.Range("A14").Value = "Household"

# Row 17
$ws.Range("B17,C17,D17,E17,F17,G17,H17").NumberFormat = "@"
$ws.Range("A17").Value = "Less than 15"
$ws.Range("B17").Value = "0.09"
$ws.Range("C17").Value = "0.27"
$ws.Range("D17").Value = "0.03"
$ws.Range("E17").Value = "0.05"
$ws.Range("F17").Value = "0.04"
$ws.Range("G17").Value = "0.03"
$ws.Range("H17").Value = "0.11"
$ws.Range("I17").Value = "-"
$ws.Range("J17").Value = "-"
$ws.Range("K17").Value = "-"
$ws.Range("L17").Value = "-"

# Row 18
$ws.Range("B18,C18,D18,E18,F18,G18,H18").NumberFormat = "@"
$ws.Range("A18").Value = "15 - 24"
$ws.Range("B18").Value = "0.22"
$ws.Range("C18").Value = "0.36"
$ws.Range("D18").Value = "0.16"
$ws.Range("E18").Value = "0.15"
$ws.Range("F18").Value = "0.16"
$ws.Range("G18").Value = "0.17"
$ws.Range("H18").Value = "0.2"

# Row 19
$ws.Range("B19,C19,D19,E19,F19,G19,H19").NumberFormat = "@"
$ws.Range("A19").Value = "25 - 49"
$ws.Range("B19").Value = "0.47"
$ws.Range("C19").Value = "0.25"
$ws.Range("D19").Value = "0.5"
$ws.Range("E19").Value = "0.55"
$ws.Range("F19").Value = "0.55"
$ws.Range("G19").Value = "0.52"
$ws.Range("H19").Value = "0.45"

# Row 20
$ws.Range("A20").Value = "Age Migrated"

# Row 22
$ws.Range("B22,C22,D22,E22,F22,G22,H22").NumberFormat = "@"
$ws.Range("A22").Value = "Before 1965"
$ws.Range("B22").Value = "0.27"
$ws.Range("C22").Value = "0.61"
$ws.Range("D22").Value = "0.17"
$ws.Range("E22").Value = "0.32"
$ws.Range("F22").Value = "0.17"
$ws.Range("G22").Value = "0.18"
$ws.Range("H22").Value = "0.33"
$ws.Range("I22").Value = "-"
$ws.Range("J22").Value = "-"
$ws.Range("K22").Value = "-"
$ws.Range("L22").Value = "-"

# Row 23
$ws.Range("B23,C23,D23,E23,F23,G23,H23").NumberFormat = "@"
$ws.Range("A23").Value = "1965 - 1979"
$ws.Range("B23").Value = "0.4"
$ws.Range("C23").Value = "0.21"
$ws.Range("D23").Value = "0.36"
$ws.Range("E23").Value = "0.4"
$ws.Range("F23").Value = "0.37"
$ws.Range("G23").Value = "0.4"
$ws.Range("H23").Value = "0.32"

# Row 24
$ws.Range("B24,C24,D24,E24,F24,G24,H24").NumberFormat = "@"
$ws.Range("A24").Value = "1980 - 1999"
$ws.Range("B24").Value = "0.26"
$ws.Range("C24").Value = "0.12"
$ws.Range("D24").Value = "0.37"
$ws.Range("E24").Value = "0.27"
$ws.Range("F24").Value = "0.42"
$ws.Range("G24").Value = "0.31"
$ws.Range("H24").Value = "0.29"

# Row 25
$ws.Range("B25,C25,D25,E25,F25,G25,H25").NumberFormat = "@"
$ws.Range("A25").Value = "After 1999"
$ws.Range("B25").Value = "0.1"
$ws.Range("C25").Value = "0.07"
$ws.Range("D25").Value = "0.14"
$ws.Range("E25").Value = "0.1"
$ws.Range("F25").Value = "0.1"
$ws.Range("G25").Value = "0.14"
$ws.Range("H25").Value = "0.09"

# Row 26
$ws.Range("A26").Value = "Migration Cohort"

# Row 27
$ws.Range("B27,D27,E27,F27,G27,H27").NumberFormat = "@"
$ws.Range("A27").Value = "Citizen"
$ws.Range("B27").Value = "0.49"
$ws.Range("C27").Value = "-"
$ws.Range("D27").Value = "0.58"
$ws.Range("E27").Value = "0.77"
$ws.Range("F27").Value = "0.59"
$ws.Range("G27").Value = "0.65"
$ws.Range("H27").Value = "0.73"
$ws.Range("I27").Value = "-"
$ws.Range("J27").Value = "-"
$ws.Range("K27").Value = "-"
$ws.Range("L27").Value = "-"

# Row 28
$ws.Range("B28,C28,D28,E28,F28,G28,H28,I28,J28,K28,L28").NumberFormat = "@"
$ws.Range("A28").Value = "English Speakers"
$ws.Range("B28").Value = "0.63"
$ws.Range("C28").Value = "0.91"
$ws.Range("D28").Value = "0.62"
$ws.Range("E28").Value = "0.73"
$ws.Range("F28").Value = "0.74"
$ws.Range("G28").Value = "0.83"
$ws.Range("H28").Value = "0.91"
$ws.Range("I28").Value = "0.98"
$ws.Range("J28").Value = "1"
$ws.Range("K28").Value = "1"
$ws.Range("L28").Value = "0.99"

# Row 29
$ws.Range("B29,C29,D29,E29,F29,G29,H29,I29,J29,K29,L29").NumberFormat = "@"
$ws.Range("A29").Value = "N"
$ws.Range("B29").Value = "29857"
$ws.Range("C29").Value = "10973"
$ws.Range("D29").Value = "3116"
$ws.Range("E29").Value = "11093"
$ws.Range("F29").Value = "6702"
$ws.Range("G29").Value = "9479"
$ws.Range("H29").Value = "154527"
$ws.Range("I29").Value = "47087"
$ws.Range("J29").Value = "141695"
$ws.Range("K29").Value = "1551834"
$ws.Range("L29").Value = "34067"

# Row 30
$ws.Range("A30").Value = "Acculturation"

